# Weekly update: "Fruta / hortaliza, semanal"
#
# The sheet is a weekly price log for "Pepino ensalada" (Comercializadora
# del Agro de Limarí). A new week's worth of data (2 rows: "Primera" and
# "Segunda" quality) is inserted right before the existing row 125 block,
# pushing the previous weeks down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 123, shifting the old
# rows 123-127 down to 125-129 (their values are left untouched by the
# insert itself).
$ws.Rows("123:124").Insert()

# Fill the newly inserted row 123 with the new week's "Primera" quality data.
$ws.Range("A123").Value = 2
$ws.Range("B123").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C123").Value = "Coquimbo"
$ws.Range("D123").Value = 44595
$ws.Range("E123").Value = 4
$ws.Range("F123").Value = 100112043
$ws.Range("G123").Value = "Pepino ensalada"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 600
$ws.Range("K123").Value = 11000
$ws.Range("L123").Value = 12000
$ws.Range("M123").Value = 11500
$ws.Range("N123").Value = "$/caja 70 unidades"
$ws.Range("O123").Value = "Provincia de Limarí"
$ws.Range("P123").Value = 164
$ws.Range("Q123").Value = 70
$ws.Range("R123").Value = "Hortaliza"

# Fill the newly inserted row 124 with the new week's "Segunda" quality data.
$ws.Range("A124").Value = 2
$ws.Range("B124").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C124").Value = "Coquimbo"
$ws.Range("D124").Value = 44595
$ws.Range("E124").Value = 4
$ws.Range("F124").Value = 100112043
$ws.Range("G124").Value = "Pepino ensalada"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Segunda"
$ws.Range("J124").Value = 400
$ws.Range("K124").Value = 8000
$ws.Range("L124").Value = 9000
$ws.Range("M124").Value = 8500
$ws.Range("N124").Value = "$/caja 100 unidades"
$ws.Range("O124").Value = "Provincia de Limarí"
$ws.Range("P124").Value = 85
$ws.Range("Q124").Value = 100
$ws.Range("R124").Value = "Hortaliza"
